$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.932.58'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.431.15'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.94'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.61'
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +6.34%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.736'
$ws.Range('E9').Value = '  +7.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.139'
$ws.Range('E10').Value = '  +6.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.72'
$ws.Range('E11').Value = '  +2.81%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.65'
$ws.Range('E12').Value = '  +9.47%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.14'
$ws.Range('E13').Value = '  +9.49%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.964.32'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('E16').Value = '  +40.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.441.58'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.29'
$ws.Range('E18').Value = '  +5.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.08'
$ws.Range('E19').Value = '  +6.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '61.868.30'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '445.53'
$ws.Range('E21').Value = '  +43.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '92.12'
$ws.Range('E22').Value = '  +10.53%  '
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.99'
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.25'
$ws.Range('E25').Value = '  +3.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '32.86'
$ws.Range('E26').Value = '  +11.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.69'
$ws.Range('E27').Value = '  +7.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.76'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.70'
$ws.Range('E29').Value = '  -3.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.73'
$ws.Range('E30').Value = '  +1.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.00'
$ws.Range('E31').Value = '  +6.45%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.03'
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0499'
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.23'
$ws.Range('E37').Value = '  +3.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('E40').Value = '  +8.19%  '
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.318'
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '142.00'
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.22'
$ws.Range('E44').Value = '  +7.33%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.54'
$ws.Range('E45').Value = '  +14.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.99'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.61'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.46'
$ws.Range('E48').Value = '  +5.99%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.13'
$ws.Range('E49').Value = '  +12.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.139'
$ws.Range('E50').Value = '  +20.21%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.773.70'
$ws.Range('E51').Value = '  +1.09%  '
